$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing _GoBack bookmark; we'll re-create it at the
#    new end of the document once all the new content is in place.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Append two new "NoSpacing" paragraphs after the current last
#    paragraph.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tail = $lastPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$tail = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

# ------------------------------------------------------------------
# 3. Fill the first new paragraph:
#    <br/>I don't think I like interacting with an alignment through
#    an oracle, actually.  + I'd rather have a function ...
# ------------------------------------------------------------------
$para6 = $d.Paragraphs.Item(6).Range
$para6.Collapse(1)
$para6.InsertAfter("@@BR@@I don't think I like interacting with an alignment through an oracle, actually. I'd rather have a function of a seq record, and just use the oracle to pass along the alignment along with info about which is the pdb sequence. So, a class that inherits from the alignment class.")

$f = $d.Paragraphs.Item(6).Range.Find
$f.Execute("@@BR@@", $false, $false, $false, $false, $false, $true, 1, $false, "^l", 2)

# ------------------------------------------------------------------
# 4. Fill the second new paragraph:
#    Thing is I find myself ... object just [italic]was[/italic] an
#    alignment.
# ------------------------------------------------------------------
$para7 = $d.Paragraphs.Item(7).Range
$para7.Collapse(1)
$para7.InsertAfter("Thing is I find myself just pulling sequences out of the alignment to use the compare function. It would be better if the object just was an alignment.")

$target = $d.Paragraphs.Item(7).Range
$target.Find.Execute("was", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target.Italic = 1

# ------------------------------------------------------------------
# 5. Re-create the _GoBack bookmark, collapsed, right after the very
#    last character of the document content (matches original
#    placement relative to the last run of text).
# ------------------------------------------------------------------
$endPos = $d.Content.End - 1
$marker = $d.Range($endPos, $endPos)
$marker.InsertAfter("@@END@@")

$bmPos = $d.Content.End - 8
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanupPos = $d.Content.End - 8
$cleanupRange = $d.Range($cleanupPos, $cleanupPos + 7)
$cleanupRange.Delete()
